# Update cryptocurrency price (D) and volume-change (E) columns
# for rows 2-51 on the active sheet, per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.907.00'
$ws.Range("E2").Value = '  -5.06%  '

$ws.Range("D3").Value = '2.468.25'
$ws.Range("E3").Value = '  -8.30%  '

$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.29%  '

$ws.Range("D5").Value = "'469.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.44%  '

$ws.Range("D6").Value = "'134.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.35%  '

$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("D8").Value = "'0.492"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.20%  '

$ws.Range("D9").Value = '2.477.76'
$ws.Range("E9").Value = '  -8.18%  '

$ws.Range("D10").Value = "'0.0966"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.71%  '

$ws.Range("E11").Value = '  -8.57%  '

$ws.Range("E12").Value = '  -5.97%  '

$ws.Range("E13").Value = '  -3.30%  '

$ws.Range("D14").Value = '2.881.08'
$ws.Range("E14").Value = '  -8.97%  '

$ws.Range("D15").Value = '54.629.59'
$ws.Range("E15").Value = '  -5.72%  '

$ws.Range("D16").Value = "'20.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.80%  '

$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("D18").Value = '2.477.72'
$ws.Range("E18").Value = '  -8.24%  '

$ws.Range("E19").Value = '  -8.40%  '

$ws.Range("D20").Value = "'311.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.34%  '

$ws.Range("D21").Value = "'9.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -11.39%  '

$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").Value = "'5.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.52%  '

$ws.Range("D24").Value = "'5.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.68%  '

$ws.Range("D25").Value = "'57.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.96%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("D27").Value = "'0.388"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.12%  '

$ws.Range("D28").Value = '2.552.68'
$ws.Range("E28").Value = '  -9.94%  '

$ws.Range("D29").Value = "'0.156"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.76%  '

$ws.Range("D30").Value = "'7.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").Value = "'0.996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.22%  '

$ws.Range("D32").Value = '0.0₃0733'
$ws.Range("E32").Value = '  -9.04%  '

$ws.Range("D33").Value = "'149.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.93%  '

$ws.Range("D34").Value = "'17.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.37%  '

$ws.Range("E35").Value = '  -8.19%  '

$ws.Range("E36").Value = '  -2.78%  '

$ws.Range("E37").Value = '  -12.61%  '

$ws.Range("E38").Value = '  -3.23%  '

$ws.Range("D39").Value = "'0.809"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.99%  '

$ws.Range("D40").Value = "'33.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.80%  '

$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("E42").Value = '  +3.33%  '

$ws.Range("D43").Value = "'0.0533"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '

$ws.Range("D44").Value = "'3.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.42%  '

$ws.Range("E45").Value = '  -5.96%  '

$ws.Range("D46").Value = "'10.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("D47").Value = '1.965.28'
$ws.Range("E47").Value = '  -8.09%  '

$ws.Range("D48").Value = "'0.0220"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.92%  '

$ws.Range("D49").Value = "'0.0881"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.37%  '

$ws.Range("D50").Value = "'4.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.17%  '

$ws.Range("D51").Value = "'16.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.67%  '
